$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B -> C, old C -> D)
$ws.Range("B1").EntireColumn.Insert()

# Header row: new StatQuery column header
$ws.Range("B1").Value = "StatQuery"

# New column B takes on the same width as column A
$ws.Range("B1").EntireColumn.ColumnWidth = $ws.Range("A1").EntireColumn.ColumnWidth

# Row 2: the two Cypher query strings (new A2/B2 content)
$ws.Range("A2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s WHERE f.file_type IN [''Aligned DNA reads file'']  RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'
$ws.Range("B2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report)OPTIONAL MATCH (s)<-[*]-(f:file)  WITH  c AS c, t ,a, s , f WHERE f.file_type IN [''Aligned DNA reads file'',''Aligned RNA reads file'',''Index file'',''Variants file''] RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trials'

# Wrap text + taller row to show the multi-line queries
$ws.Range("A2:B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 101.5

# Move selection, matching the saved view state
$ws.Range("A4").Select()

$wb.Save()
